# chore: update Sheets via scheduled runner
#
# Refreshes the cached market-board price columns (H:N -- currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) across all eight
# crafter-job leve tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with newly
# pulled values. A handful of rows gain or lose a trailing N (LeveProfitHQ)
# or M (LeveProfitNQ) cell depending on whether an HQ/NQ price is available.

$wb = $excel.ActiveWorkbook


# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(8, 8).Value = 74.5  # H8: 87.22221999999999 -> 74.5
$ws.Cells.Item(8, 9).Value = 74.5  # I8: 83.71429000000001 -> 74.5
$ws.Cells.Item(8, 10).Value = 0  # J8: 99.5 -> 0
$ws.Cells.Item(8, 11).Value = 223.5  # K8: 251.14287 -> 223.5
$ws.Cells.Item(8, 12).Value = 0  # L8: 298.5 -> 0
$ws.Cells.Item(8, 13).Value = -84.5  # M8: -112.14287 -> -84.5
$ws.Cells.Item(8, 14).Value = $null  # N8: -576.5 -> None
$ws.Cells.Item(17, 8).Value = 200  # H17: 235.6875 -> 200
$ws.Cells.Item(17, 10).Value = 200  # J17: 235.6875 -> 200
$ws.Cells.Item(17, 12).Value = 600  # L17: 707.0625 -> 600
$ws.Cells.Item(17, 14).Value = -936  # N17: -1043.0625 -> -936
$ws.Cells.Item(32, 8).Value = 4205.625  # H32: 4441.8184 -> 4205.625
$ws.Cells.Item(32, 10).Value = 3654.5  # J32: 4051.9048 -> 3654.5
$ws.Cells.Item(32, 12).Value = 3654.5  # L32: 4051.9048 -> 3654.5
$ws.Cells.Item(32, 14).Value = -4306.5  # N32: -4703.9048 -> -4306.5
$ws.Cells.Item(40, 8).Value = 3530.6155  # H40: 3921.1428 -> 3530.6155
$ws.Cells.Item(40, 10).Value = 9666.333000000001  # J40: 9499.25 -> 9666.333000000001
$ws.Cells.Item(40, 12).Value = 9666.333000000001  # L40: 9499.25 -> 9666.333000000001
$ws.Cells.Item(40, 14).Value = -10016.333  # N40: -9849.25 -> -10016.333
$ws.Cells.Item(57, 8).Value = 49989  # H57: 69994 -> 49989
$ws.Cells.Item(57, 10).Value = 49989  # J57: 69994 -> 49989
$ws.Cells.Item(57, 12).Value = 149967  # L57: 209982 -> 149967
$ws.Cells.Item(57, 14).Value = -150965  # N57: -210980 -> -150965
$ws.Cells.Item(62, 8).Value = 0  # H62: 1499.5 -> 0
$ws.Cells.Item(62, 9).Value = 0  # I62: 1499.5 -> 0
$ws.Cells.Item(62, 11).Value = 0  # K62: 1499.5 -> 0
$ws.Cells.Item(62, 13).Value = $null  # M62: -875.5 -> None
$ws.Cells.Item(65, 8).Value = 0  # H65: 1499.5 -> 0
$ws.Cells.Item(65, 9).Value = 0  # I65: 1499.5 -> 0
$ws.Cells.Item(65, 11).Value = 0  # K65: 7497.5 -> 0
$ws.Cells.Item(65, 13).Value = $null  # M65: -4377.5 -> None
$ws.Cells.Item(80, 8).Value = 927.2632  # H80: 927.5294 -> 927.2632
$ws.Cells.Item(80, 9).Value = 693.25  # I80: 573.75 -> 693.25
$ws.Cells.Item(80, 10).Value = 1097.4546  # J80: 1242 -> 1097.4546
$ws.Cells.Item(80, 11).Value = 2079.75  # K80: 1721.25 -> 2079.75
$ws.Cells.Item(80, 12).Value = 3292.3638  # L80: 3726 -> 3292.3638
$ws.Cells.Item(80, 13).Value = -1081.75  # M80: -723.25 -> -1081.75
$ws.Cells.Item(80, 14).Value = -5288.3638  # N80: -5722 -> -5288.3638
$ws.Cells.Item(83, 8).Value = 927.2632  # H83: 927.5294 -> 927.2632
$ws.Cells.Item(83, 9).Value = 693.25  # I83: 573.75 -> 693.25
$ws.Cells.Item(83, 10).Value = 1097.4546  # J83: 1242 -> 1097.4546
$ws.Cells.Item(83, 11).Value = 6239.25  # K83: 5163.75 -> 6239.25
$ws.Cells.Item(83, 12).Value = 9877.091400000001  # L83: 11178 -> 9877.091400000001
$ws.Cells.Item(83, 13).Value = -1247.25  # M83: -171.75 -> -1247.25
$ws.Cells.Item(83, 14).Value = -19861.0914  # N83: -21162 -> -19861.0914
$ws.Cells.Item(86, 8).Value = 5250  # H86: 5071.4287 -> 5250
$ws.Cells.Item(86, 9).Value = 4900  # I86: 4750 -> 4900
$ws.Cells.Item(86, 11).Value = 4900  # K86: 4750 -> 4900
$ws.Cells.Item(86, 13).Value = -3777  # M86: -3627 -> -3777
$ws.Cells.Item(89, 8).Value = 5250  # H89: 5071.4287 -> 5250
$ws.Cells.Item(89, 9).Value = 4900  # I89: 4750 -> 4900
$ws.Cells.Item(89, 11).Value = 24500  # K89: 23750 -> 24500
$ws.Cells.Item(89, 13).Value = -18884  # M89: -18134 -> -18884
$ws.Cells.Item(92, 8).Value = 30940.633  # H92: 28254.637 -> 30940.633
$ws.Cells.Item(92, 9).Value = 723.3333  # I92: 797.9259 -> 723.3333
$ws.Cells.Item(92, 11).Value = 723.3333  # K92: 797.9259 -> 723.3333
$ws.Cells.Item(92, 13).Value = 524.6667  # M92: 450.0741 -> 524.6667
$ws.Cells.Item(106, 8).Value = 8321.333000000001  # H106: 6949.4 -> 8321.333000000001
$ws.Cells.Item(106, 10).Value = 0  # J106: 4891.5 -> 0
$ws.Cells.Item(106, 12).Value = 0  # L106: 4891.5 -> 0
$ws.Cells.Item(106, 14).Value = $null  # N106: -6153.5 -> None
$ws.Cells.Item(111, 8).Value = 2914.5  # H111: 3000 -> 2914.5
$ws.Cells.Item(111, 9).Value = 2914.5  # I111: 3000 -> 2914.5
$ws.Cells.Item(111, 11).Value = 8743.5  # K111: 9000 -> 8743.5
$ws.Cells.Item(111, 13).Value = -5676.5  # M111: -5933 -> -5676.5
$ws.Cells.Item(112, 8).Value = 1731.8889  # H112: 1760.5769 -> 1731.8889
$ws.Cells.Item(112, 10).Value = 1948.55  # J112: 1999.2106 -> 1948.55
$ws.Cells.Item(112, 12).Value = 5845.65  # L112: 5997.6318 -> 5845.65
$ws.Cells.Item(112, 14).Value = -8061.65  # N112: -8213.631799999999 -> -8061.65
$ws.Cells.Item(115, 8).Value = 2436.3572  # H115: 294.26666 -> 2436.3572
$ws.Cells.Item(115, 9).Value = 300.3846  # I115: 294.26666 -> 300.3846
$ws.Cells.Item(115, 10).Value = 30204  # J115: 0 -> 30204
$ws.Cells.Item(115, 11).Value = 901.1537999999999  # K115: 882.79998 -> 901.1537999999999
$ws.Cells.Item(115, 12).Value = 90612  # L115: 0 -> 90612
$ws.Cells.Item(115, 13).Value = 665.8462000000001  # M115: 684.20002 -> 665.8462000000001
$ws.Cells.Item(115, 14).Value = -93746  # N115: None -> -93746
$ws.Cells.Item(116, 8).Value = 7854.927  # H116: 7898.737 -> 7854.927
$ws.Cells.Item(116, 9).Value = 7468.6665  # I116: 7487.407 -> 7468.6665
$ws.Cells.Item(116, 11).Value = 7468.6665  # K116: 7487.407 -> 7468.6665
$ws.Cells.Item(116, 13).Value = -4026.6665  # M116: -4045.407 -> -4026.6665
$ws.Cells.Item(118, 8).Value = 1147.5  # H118: 1333 -> 1147.5
$ws.Cells.Item(118, 9).Value = 863.3333  # I118: 999.5 -> 863.3333
$ws.Cells.Item(118, 11).Value = 2589.9999  # K118: 2998.5 -> 2589.9999
$ws.Cells.Item(118, 13).Value = -932.9998999999998  # M118: -1341.5 -> -932.9998999999998
$ws.Cells.Item(132, 8).Value = 1044563.75  # H132: 1114079.8 -> 1044563.75
$ws.Cells.Item(132, 9).Value = 1044563.75  # I132: 1114079.8 -> 1044563.75
$ws.Cells.Item(132, 11).Value = 3133691.25  # K132: 3342239.4 -> 3133691.25
$ws.Cells.Item(132, 13).Value = -3131161.25  # M132: -3339709.4 -> -3131161.25
$ws.Cells.Item(136, 8).Value = 72000  # H136: 0 -> 72000
$ws.Cells.Item(136, 10).Value = 72000  # J136: 0 -> 72000
$ws.Cells.Item(136, 12).Value = 72000  # L136: 0 -> 72000
$ws.Cells.Item(136, 14).Value = -82200  # N136: None -> -82200
$ws.Cells.Item(137, 8).Value = 1515.9445  # H137: 1426 -> 1515.9445
$ws.Cells.Item(137, 9).Value = 1380.5  # I137: 1337.2778 -> 1380.5
$ws.Cells.Item(137, 10).Value = 2599.5  # J137: 2224.5 -> 2599.5
$ws.Cells.Item(137, 11).Value = 4141.5  # K137: 4011.8334 -> 4141.5
$ws.Cells.Item(137, 12).Value = 7798.5  # L137: 6673.5 -> 7798.5
$ws.Cells.Item(137, 13).Value = -1591.5  # M137: -1461.8334 -> -1591.5
$ws.Cells.Item(137, 14).Value = -12898.5  # N137: -11773.5 -> -12898.5
$ws.Cells.Item(139, 8).Value = 100000  # H139: 0 -> 100000
$ws.Cells.Item(139, 10).Value = 100000  # J139: 0 -> 100000
$ws.Cells.Item(139, 12).Value = 100000  # L139: 0 -> 100000
$ws.Cells.Item(139, 14).Value = -110280  # N139: None -> -110280

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1760.35  # H2: 1760.9 -> 1760.35
$ws.Cells.Item(2, 9).Value = 1856.2  # I2: 1856.2667 -> 1856.2
$ws.Cells.Item(2, 10).Value = 1472.8  # J2: 1474.8 -> 1472.8
$ws.Cells.Item(2, 11).Value = 1856.2  # K2: 1856.2667 -> 1856.2
$ws.Cells.Item(2, 12).Value = 1472.8  # L2: 1474.8 -> 1472.8
$ws.Cells.Item(2, 13).Value = -1743.2  # M2: -1743.2667 -> -1743.2
$ws.Cells.Item(2, 14).Value = -1698.8  # N2: -1700.8 -> -1698.8
$ws.Cells.Item(32, 8).Value = 3201.5874  # H32: 3011.5671 -> 3201.5874
$ws.Cells.Item(32, 9).Value = 2835.1133  # I32: 2637.4736 -> 2835.1133
$ws.Cells.Item(32, 11).Value = 2835.1133  # K32: 2637.4736 -> 2835.1133
$ws.Cells.Item(32, 13).Value = -2548.1133  # M32: -2350.4736 -> -2548.1133
$ws.Cells.Item(34, 8).Value = 38749.832  # H34: 36500 -> 38749.832
$ws.Cells.Item(34, 9).Value = 37499  # I34: 24000 -> 37499
$ws.Cells.Item(34, 11).Value = 37499  # K34: 24000 -> 37499
$ws.Cells.Item(34, 13).Value = -37228  # M34: -23729 -> -37228
$ws.Cells.Item(61, 8).Value = 4367.7617  # H61: 4486.25 -> 4367.7617
$ws.Cells.Item(61, 9).Value = 3147.3333  # I61: 3291 -> 3147.3333
$ws.Cells.Item(61, 11).Value = 3147.3333  # K61: 3291 -> 3147.3333
$ws.Cells.Item(61, 13).Value = -2935.3333  # M61: -3079 -> -2935.3333
$ws.Cells.Item(74, 8).Value = 2503.7908  # H74: 3099.12 -> 2503.7908
$ws.Cells.Item(74, 9).Value = 1886.129  # I74: 2083.7856 -> 1886.129
$ws.Cells.Item(74, 10).Value = 4099.4165  # J74: 4391.364 -> 4099.4165
$ws.Cells.Item(74, 11).Value = 1886.129  # K74: 2083.7856 -> 1886.129
$ws.Cells.Item(74, 12).Value = 4099.4165  # L74: 4391.364 -> 4099.4165
$ws.Cells.Item(74, 13).Value = -1012.129  # M74: -1209.7856 -> -1012.129
$ws.Cells.Item(74, 14).Value = -5847.4165  # N74: -6139.364 -> -5847.4165
$ws.Cells.Item(77, 8).Value = 2503.7908  # H77: 3099.12 -> 2503.7908
$ws.Cells.Item(77, 9).Value = 1886.129  # I77: 2083.7856 -> 1886.129
$ws.Cells.Item(77, 10).Value = 4099.4165  # J77: 4391.364 -> 4099.4165
$ws.Cells.Item(77, 11).Value = 9430.645  # K77: 10418.928 -> 9430.645
$ws.Cells.Item(77, 12).Value = 20497.0825  # L77: 21956.82 -> 20497.0825
$ws.Cells.Item(77, 13).Value = -5062.645  # M77: -6050.928 -> -5062.645
$ws.Cells.Item(77, 14).Value = -29233.0825  # N77: -30692.82 -> -29233.0825
$ws.Cells.Item(116, 8).Value = 1760.35  # H116: 1760.9 -> 1760.35
$ws.Cells.Item(116, 9).Value = 1856.2  # I116: 1856.2667 -> 1856.2
$ws.Cells.Item(116, 10).Value = 1472.8  # J116: 1474.8 -> 1472.8
$ws.Cells.Item(116, 11).Value = 1856.2  # K116: 1856.2667 -> 1856.2
$ws.Cells.Item(116, 12).Value = 1472.8  # L116: 1474.8 -> 1472.8
$ws.Cells.Item(116, 13).Value = 437.8  # M116: 437.7333000000001 -> 437.8
$ws.Cells.Item(116, 14).Value = -6060.8  # N116: -6062.8 -> -6060.8
$ws.Cells.Item(132, 8).Value = 5130969.5  # H132: 5265975.5 -> 5130969.5
$ws.Cells.Item(132, 9).Value = 2009.4138  # I132: 2054.4285 -> 2009.4138
$ws.Cells.Item(132, 11).Value = 6028.2414  # K132: 6163.2855 -> 6028.2414
$ws.Cells.Item(132, 13).Value = -3498.2414  # M132: -3633.2855 -> -3498.2414
$ws.Cells.Item(136, 8).Value = 4367.7617  # H136: 4486.25 -> 4367.7617
$ws.Cells.Item(136, 9).Value = 3147.3333  # I136: 3291 -> 3147.3333
$ws.Cells.Item(136, 11).Value = 9441.999899999999  # K136: 9873 -> 9441.999899999999
$ws.Cells.Item(136, 13).Value = -6891.999899999999  # M136: -7323 -> -6891.999899999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1760.35  # H3: 1760.9 -> 1760.35
$ws.Cells.Item(3, 9).Value = 1856.2  # I3: 1856.2667 -> 1856.2
$ws.Cells.Item(3, 10).Value = 1472.8  # J3: 1474.8 -> 1472.8
$ws.Cells.Item(3, 11).Value = 1856.2  # K3: 1856.2667 -> 1856.2
$ws.Cells.Item(3, 12).Value = 1472.8  # L3: 1474.8 -> 1472.8
$ws.Cells.Item(3, 13).Value = -1742.2  # M3: -1742.2667 -> -1742.2
$ws.Cells.Item(3, 14).Value = -1700.8  # N3: -1702.8 -> -1700.8
$ws.Cells.Item(22, 8).Value = 714.8946999999999  # H22: 727.3889 -> 714.8946999999999
$ws.Cells.Item(22, 9).Value = 681.41174  # I22: 693.375 -> 681.41174
$ws.Cells.Item(22, 11).Value = 681.41174  # K22: 693.375 -> 681.41174
$ws.Cells.Item(22, 13).Value = -508.41174  # M22: -520.375 -> -508.41174
$ws.Cells.Item(74, 8).Value = 53980.715  # H74: 55979.332 -> 53980.715
$ws.Cells.Item(74, 9).Value = 41989  # I74: 0 -> 41989
$ws.Cells.Item(74, 11).Value = 41989  # K74: 0 -> 41989
$ws.Cells.Item(74, 13).Value = -41053  # M74: None -> -41053
$ws.Cells.Item(77, 8).Value = 53980.715  # H77: 55979.332 -> 53980.715
$ws.Cells.Item(77, 9).Value = 41989  # I77: 0 -> 41989
$ws.Cells.Item(77, 11).Value = 125967  # K77: 0 -> 125967
$ws.Cells.Item(77, 13).Value = -121287  # M77: None -> -121287
$ws.Cells.Item(80, 8).Value = 423.7143  # H80: 451.85715 -> 423.7143
$ws.Cells.Item(80, 10).Value = 436.2  # J80: 475.6 -> 436.2
$ws.Cells.Item(80, 12).Value = 436.2  # L80: 475.6 -> 436.2
$ws.Cells.Item(80, 14).Value = -2432.2  # N80: -2471.6 -> -2432.2
$ws.Cells.Item(83, 8).Value = 423.7143  # H83: 451.85715 -> 423.7143
$ws.Cells.Item(83, 10).Value = 436.2  # J83: 475.6 -> 436.2
$ws.Cells.Item(83, 12).Value = 2181  # L83: 2378 -> 2181
$ws.Cells.Item(83, 14).Value = -12165  # N83: -12362 -> -12165
$ws.Cells.Item(86, 8).Value = 38462456  # H86: 41667600 -> 38462456
$ws.Cells.Item(86, 9).Value = 66667540  # I86: 76923976 -> 66667540
$ws.Cells.Item(86, 11).Value = 66667540  # K86: 76923976 -> 66667540
$ws.Cells.Item(86, 13).Value = -66666417  # M86: -76922853 -> -66666417
$ws.Cells.Item(89, 8).Value = 38462456  # H89: 41667600 -> 38462456
$ws.Cells.Item(89, 9).Value = 66667540  # I89: 76923976 -> 66667540
$ws.Cells.Item(89, 11).Value = 333337700  # K89: 384619880 -> 333337700
$ws.Cells.Item(89, 13).Value = -333332084  # M89: -384614264 -> -333332084
$ws.Cells.Item(107, 8).Value = 42282.844  # H107: 42285.21 -> 42282.844
$ws.Cells.Item(107, 9).Value = 56098.355  # I107: 56098.43 -> 56098.355
$ws.Cells.Item(107, 10).Value = 3599.4  # J107: 3608.2 -> 3599.4
$ws.Cells.Item(107, 11).Value = 56098.355  # K107: 56098.43 -> 56098.355
$ws.Cells.Item(107, 12).Value = 3599.4  # L107: 3608.2 -> 3599.4
$ws.Cells.Item(107, 13).Value = -54178.355  # M107: -54178.43 -> -54178.355
$ws.Cells.Item(107, 14).Value = -7439.4  # N107: -7448.2 -> -7439.4

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 277  # H7: 279.94116 -> 277
$ws.Cells.Item(7, 10).Value = 468.66666  # J7: 477 -> 468.66666
$ws.Cells.Item(7, 12).Value = 468.66666  # L7: 477 -> 468.66666
$ws.Cells.Item(7, 14).Value = -694.66666  # N7: -703 -> -694.66666
$ws.Cells.Item(16, 8).Value = 1483.2222  # H16: 1591.1428 -> 1483.2222
$ws.Cells.Item(16, 9).Value = 1163  # I16: 1220.5 -> 1163
$ws.Cells.Item(16, 11).Value = 1163  # K16: 1220.5 -> 1163
$ws.Cells.Item(16, 13).Value = -876  # M16: -933.5 -> -876
$ws.Cells.Item(62, 8).Value = 171810.83  # H62: 105732.6 -> 171810.83
$ws.Cells.Item(62, 9).Value = 171810.83  # I62: 105732.6 -> 171810.83
$ws.Cells.Item(62, 11).Value = 171810.83  # K62: 105732.6 -> 171810.83
$ws.Cells.Item(62, 13).Value = -171186.83  # M62: -105108.6 -> -171186.83
$ws.Cells.Item(64, 8).Value = 46332.133  # H64: 49000 -> 46332.133
$ws.Cells.Item(64, 10).Value = 46332.133  # J64: 49000 -> 46332.133
$ws.Cells.Item(64, 12).Value = 46332.133  # L64: 49000 -> 46332.133
$ws.Cells.Item(64, 14).Value = -46828.133  # N64: -49496 -> -46828.133
$ws.Cells.Item(65, 8).Value = 171810.83  # H65: 105732.6 -> 171810.83
$ws.Cells.Item(65, 9).Value = 171810.83  # I65: 105732.6 -> 171810.83
$ws.Cells.Item(65, 11).Value = 859054.1499999999  # K65: 528663 -> 859054.1499999999
$ws.Cells.Item(65, 13).Value = -855934.1499999999  # M65: -525543 -> -855934.1499999999
$ws.Cells.Item(67, 8).Value = 46332.133  # H67: 49000 -> 46332.133
$ws.Cells.Item(67, 10).Value = 46332.133  # J67: 49000 -> 46332.133
$ws.Cells.Item(67, 12).Value = 46332.133  # L67: 49000 -> 46332.133
$ws.Cells.Item(67, 14).Value = -48048.133  # N67: -50716 -> -48048.133
$ws.Cells.Item(107, 8).Value = 1250.1212  # H107: 1242.9429 -> 1250.1212
$ws.Cells.Item(107, 9).Value = 1016.46155  # I107: 1015.85187 -> 1016.46155
$ws.Cells.Item(107, 10).Value = 2118  # J107: 2009.375 -> 2118
$ws.Cells.Item(107, 11).Value = 1016.46155  # K107: 1015.85187 -> 1016.46155
$ws.Cells.Item(107, 12).Value = 2118  # L107: 2009.375 -> 2118
$ws.Cells.Item(107, 13).Value = 903.53845  # M107: 904.14813 -> 903.53845
$ws.Cells.Item(107, 14).Value = -5958  # N107: -5849.375 -> -5958
$ws.Cells.Item(113, 8).Value = 1483.2222  # H113: 1591.1428 -> 1483.2222
$ws.Cells.Item(113, 9).Value = 1163  # I113: 1220.5 -> 1163
$ws.Cells.Item(113, 11).Value = 1163  # K113: 1220.5 -> 1163
$ws.Cells.Item(113, 13).Value = 1007  # M113: 949.5 -> 1007
$ws.Cells.Item(132, 8).Value = 15589.2  # H132: 13616 -> 15589.2
$ws.Cells.Item(132, 9).Value = 13236.5  # I132: 13599 -> 13236.5
$ws.Cells.Item(132, 10).Value = 25000  # J132: 13650 -> 25000
$ws.Cells.Item(132, 11).Value = 39709.5  # K132: 40797 -> 39709.5
$ws.Cells.Item(132, 12).Value = 75000  # L132: 40950 -> 75000
$ws.Cells.Item(132, 13).Value = -37179.5  # M132: -38267 -> -37179.5
$ws.Cells.Item(132, 14).Value = -80060  # N132: -46010 -> -80060
$ws.Cells.Item(134, 8).Value = 22226922  # H134: 15877673 -> 22226922
$ws.Cells.Item(134, 9).Value = 22226922  # I134: 18523118 -> 22226922
$ws.Cells.Item(134, 10).Value = 0  # J134: 5000 -> 0
$ws.Cells.Item(134, 11).Value = 66680766  # K134: 55569354 -> 66680766
$ws.Cells.Item(134, 12).Value = 0  # L134: 15000 -> 0
$ws.Cells.Item(134, 13).Value = -66678231  # M134: -55566819 -> -66678231
$ws.Cells.Item(134, 14).Value = $null  # N134: -20070 -> None
$ws.Cells.Item(139, 8).Value = 79460.836  # H139: 79953 -> 79460.836
$ws.Cells.Item(139, 10).Value = 79460.836  # J139: 79953 -> 79460.836
$ws.Cells.Item(139, 12).Value = 79460.836  # L139: 79953 -> 79460.836
$ws.Cells.Item(139, 14).Value = -89740.836  # N139: -90233 -> -89740.836

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1586.4166  # H5: 1612.0834 -> 1586.4166
$ws.Cells.Item(5, 9).Value = 586.46155  # I5: 619.1429000000001 -> 586.46155
$ws.Cells.Item(5, 10).Value = 2768.182  # J5: 3002.2 -> 2768.182
$ws.Cells.Item(5, 11).Value = 1759.38465  # K5: 1857.4287 -> 1759.38465
$ws.Cells.Item(5, 12).Value = 8304.545999999998  # L5: 9006.599999999999 -> 8304.545999999998
$ws.Cells.Item(5, 13).Value = -1647.38465  # M5: -1745.4287 -> -1647.38465
$ws.Cells.Item(5, 14).Value = -8528.545999999998  # N5: -9230.599999999999 -> -8528.545999999998
$ws.Cells.Item(23, 8).Value = 150.6  # H23: 151.1 -> 150.6
$ws.Cells.Item(23, 10).Value = 248.5  # J23: 249.75 -> 248.5
$ws.Cells.Item(23, 12).Value = 745.5  # L23: 749.25 -> 745.5
$ws.Cells.Item(23, 14).Value = -1215.5  # N23: -1219.25 -> -1215.5
$ws.Cells.Item(51, 8).Value = 9943.166999999999  # H51: 10931.8 -> 9943.166999999999
$ws.Cells.Item(51, 10).Value = 12189.375  # J51: 14585.833 -> 12189.375
$ws.Cells.Item(51, 12).Value = 36568.125  # L51: 43757.499 -> 36568.125
$ws.Cells.Item(51, 14).Value = -37488.125  # N51: -44677.499 -> -37488.125
$ws.Cells.Item(96, 8).Value = 12903.8  # H96: 12925.25 -> 12903.8
$ws.Cells.Item(96, 10).Value = 12903.8  # J96: 12925.25 -> 12903.8
$ws.Cells.Item(96, 12).Value = 38711.39999999999  # L96: 38775.75 -> 38711.39999999999
$ws.Cells.Item(96, 14).Value = -42829.39999999999  # N96: -42893.75 -> -42829.39999999999
$ws.Cells.Item(103, 8).Value = 2786  # H103: 2542.5 -> 2786
$ws.Cells.Item(103, 9).Value = 3240  # I103: 2896.6667 -> 3240
$ws.Cells.Item(103, 10).Value = 970  # J103: 1480 -> 970
$ws.Cells.Item(103, 11).Value = 9720  # K103: 8690.000100000001 -> 9720
$ws.Cells.Item(103, 12).Value = 2910  # L103: 4440 -> 2910
$ws.Cells.Item(103, 13).Value = -8841  # M103: -7811.000100000001 -> -8841
$ws.Cells.Item(103, 14).Value = -4668  # N103: -6198 -> -4668
$ws.Cells.Item(107, 8).Value = 545.85187  # H107: 729.9474 -> 545.85187
$ws.Cells.Item(107, 9).Value = 392.53845  # I107: 628.4286 -> 392.53845
$ws.Cells.Item(107, 10).Value = 688.2143  # J107: 789.1667 -> 688.2143
$ws.Cells.Item(107, 11).Value = 1177.61535  # K107: 1885.2858 -> 1177.61535
$ws.Cells.Item(107, 12).Value = 2064.6429  # L107: 2367.5001 -> 2064.6429
$ws.Cells.Item(107, 13).Value = 742.38465  # M107: 34.71420000000012 -> 742.38465
$ws.Cells.Item(107, 14).Value = -5904.6429  # N107: -6207.5001 -> -5904.6429
$ws.Cells.Item(121, 8).Value = 8337499.5  # H121: 9095445 -> 8337499.5
$ws.Cells.Item(121, 10).Value = 14287057  # J121: 16668216 -> 14287057
$ws.Cells.Item(121, 12).Value = 42861171  # L121: 50004648 -> 42861171
$ws.Cells.Item(121, 14).Value = -42863791  # N121: -50007268 -> -42863791
$ws.Cells.Item(135, 8).Value = 1586.4166  # H135: 1612.0834 -> 1586.4166
$ws.Cells.Item(135, 9).Value = 586.46155  # I135: 619.1429000000001 -> 586.46155
$ws.Cells.Item(135, 10).Value = 2768.182  # J135: 3002.2 -> 2768.182
$ws.Cells.Item(135, 11).Value = 5278.15395  # K135: 5572.2861 -> 5278.15395
$ws.Cells.Item(135, 12).Value = 24913.638  # L135: 27019.8 -> 24913.638
$ws.Cells.Item(135, 13).Value = -2743.15395  # M135: -3037.2861 -> -2743.15395
$ws.Cells.Item(135, 14).Value = -29983.638  # N135: -32089.8 -> -29983.638
$ws.Cells.Item(139, 8).Value = 1723.2222  # H139: 2332.7058 -> 1723.2222
$ws.Cells.Item(139, 9).Value = 1011.0833  # I139: 1787.1 -> 1011.0833
$ws.Cells.Item(139, 10).Value = 3147.5  # J139: 3112.1428 -> 3147.5
$ws.Cells.Item(139, 11).Value = 3033.2499  # K139: 5361.299999999999 -> 3033.2499
$ws.Cells.Item(139, 12).Value = 9442.5  # L139: 9336.428400000001 -> 9442.5
$ws.Cells.Item(139, 13).Value = 2106.7501  # M139: -221.2999999999993 -> 2106.7501
$ws.Cells.Item(139, 14).Value = -19722.5  # N139: -19616.4284 -> -19722.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1280.7646  # H97: 1245.6 -> 1280.7646
$ws.Cells.Item(97, 9).Value = 1147.1072  # I97: 1109.2759 -> 1147.1072
$ws.Cells.Item(97, 11).Value = 1147.1072  # K97: 1109.2759 -> 1147.1072
$ws.Cells.Item(97, 13).Value = -651.1071999999999  # M97: -613.2759000000001 -> -651.1071999999999
$ws.Cells.Item(113, 8).Value = 14257.8  # H113: 16682 -> 14257.8
$ws.Cells.Item(113, 9).Value = 16572.25  # I113: 20576 -> 16572.25
$ws.Cells.Item(113, 11).Value = 16572.25  # K113: 20576 -> 16572.25
$ws.Cells.Item(113, 13).Value = -14402.25  # M113: -18406 -> -14402.25
$ws.Cells.Item(122, 8).Value = 173156.25  # H122: 122728.234 -> 173156.25
$ws.Cells.Item(122, 9).Value = 290964  # I122: 185777 -> 290964
$ws.Cells.Item(122, 10).Value = 8225.4  # J122: 7138.8335 -> 8225.4
$ws.Cells.Item(122, 11).Value = 872892  # K122: 557331 -> 872892
$ws.Cells.Item(122, 12).Value = 24676.2  # L122: 21416.5005 -> 24676.2
$ws.Cells.Item(122, 13).Value = -870442  # M122: -554881 -> -870442
$ws.Cells.Item(122, 14).Value = -29576.2  # N122: -26316.5005 -> -29576.2

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 2633.7  # H46: 2820.5557 -> 2633.7
$ws.Cells.Item(46, 9).Value = 1113.6666  # I46: 1195.5 -> 1113.6666
$ws.Cells.Item(46, 10).Value = 3285.1428  # J46: 3284.8572 -> 3285.1428
$ws.Cells.Item(46, 11).Value = 1113.6666  # K46: 1195.5 -> 1113.6666
$ws.Cells.Item(46, 12).Value = 3285.1428  # L46: 3284.8572 -> 3285.1428
$ws.Cells.Item(46, 13).Value = -925.6666  # M46: -1007.5 -> -925.6666
$ws.Cells.Item(46, 14).Value = -3661.1428  # N46: -3660.8572 -> -3661.1428
$ws.Cells.Item(55, 8).Value = 1332.7241  # H55: 1376.7858 -> 1332.7241
$ws.Cells.Item(55, 9).Value = 426.35294  # I55: 446.8125 -> 426.35294
$ws.Cells.Item(55, 11).Value = 426.35294  # K55: 446.8125 -> 426.35294
$ws.Cells.Item(55, 13).Value = -253.35294  # M55: -273.8125 -> -253.35294
$ws.Cells.Item(61, 8).Value = 24260.1  # H61: 27650.5 -> 24260.1
$ws.Cells.Item(61, 9).Value = 22574.5  # I61: 26533.166 -> 22574.5
$ws.Cells.Item(61, 11).Value = 22574.5  # K61: 26533.166 -> 22574.5
$ws.Cells.Item(61, 13).Value = -22372.5  # M61: -26331.166 -> -22372.5
$ws.Cells.Item(113, 8).Value = 24260.1  # H113: 27650.5 -> 24260.1
$ws.Cells.Item(113, 9).Value = 22574.5  # I113: 26533.166 -> 22574.5
$ws.Cells.Item(113, 11).Value = 22574.5  # K113: 26533.166 -> 22574.5
$ws.Cells.Item(113, 13).Value = -20404.5  # M113: -24363.166 -> -20404.5
$ws.Cells.Item(122, 8).Value = 4364.641  # H122: 4395.316 -> 4364.641
$ws.Cells.Item(122, 9).Value = 3341.8147  # I122: 3347.3076 -> 3341.8147
$ws.Cells.Item(122, 11).Value = 10025.4441  # K122: 10041.9228 -> 10025.4441
$ws.Cells.Item(122, 13).Value = -7575.444100000001  # M122: -7591.9228 -> -7575.444100000001
$ws.Cells.Item(132, 8).Value = 9637.348  # H132: 10274.762 -> 9637.348
$ws.Cells.Item(132, 9).Value = 8492.333000000001  # I132: 9185.8125 -> 8492.333000000001
$ws.Cells.Item(132, 11).Value = 25476.999  # K132: 27557.4375 -> 25476.999
$ws.Cells.Item(132, 13).Value = -22946.999  # M132: -25027.4375 -> -22946.999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 4763303.5  # H81: 4203068.5 -> 4763303.5
$ws.Cells.Item(81, 9).Value = 5495935  # I81: 5953847 -> 5495935
$ws.Cells.Item(81, 10).Value = 1199  # J81: 1199.2 -> 1199
$ws.Cells.Item(81, 11).Value = 10991870  # K81: 11907694 -> 10991870
$ws.Cells.Item(81, 12).Value = 2398  # L81: 2398.4 -> 2398
$ws.Cells.Item(81, 13).Value = -10990809  # M81: -11906633 -> -10990809
$ws.Cells.Item(81, 14).Value = -4520  # N81: -4520.4 -> -4520
$ws.Cells.Item(84, 8).Value = 4763303.5  # H84: 4203068.5 -> 4763303.5
$ws.Cells.Item(84, 9).Value = 5495935  # I84: 5953847 -> 5495935
$ws.Cells.Item(84, 10).Value = 1199  # J84: 1199.2 -> 1199
$ws.Cells.Item(84, 11).Value = 54959350  # K84: 59538470 -> 54959350
$ws.Cells.Item(84, 12).Value = 11990  # L84: 11992 -> 11990
$ws.Cells.Item(84, 13).Value = -54954046  # M84: -59533166 -> -54954046
$ws.Cells.Item(84, 14).Value = -22598  # N84: -22600 -> -22598
$ws.Cells.Item(107, 8).Value = 11309.7  # H107: 11953 -> 11309.7
$ws.Cells.Item(107, 9).Value = 3073.75  # I107: 3233.0833 -> 3073.75
$ws.Cells.Item(107, 10).Value = 23663.625  # J107: 26901.428 -> 23663.625
$ws.Cells.Item(107, 11).Value = 9221.25  # K107: 9699.249899999999 -> 9221.25
$ws.Cells.Item(107, 12).Value = 70990.875  # L107: 80704.284 -> 70990.875
$ws.Cells.Item(107, 13).Value = -7301.25  # M107: -7779.249899999999 -> -7301.25
$ws.Cells.Item(107, 14).Value = -74830.875  # N107: -84544.284 -> -74830.875
